$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header strings
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# Update weekly crime statistics data cells
$data = @{
    "G15" = 2
    "H15" = 50
    "I15" = 18
    "K15" = 5.882352941176
    "L15" = 5.882352941176
    "M15" = 260
    "N15" = -30.769230769230
    "C16" = 10
    "D16" = 6
    "E16" = 66.666666666666
    "F16" = 44
    "H16" = 76
    "I16" = 179
    "J16" = 147
    "K16" = 21.768707482993
    "L16" = 31.617647058823
    "M16" = 72.115384615384
    "N16" = -61.255411255411
    "C17" = 9
    "D17" = 15
    "E17" = -40
    "F17" = 60
    "G17" = 56
    "H17" = 7.142857142857
    "I17" = 284
    "J17" = 246
    "K17" = 15.447154471544
    "L17" = 13.147410358565
    "M17" = 167.924528301887
    "N17" = -24.064171122994
    "C18" = 6
    "D18" = 3
    "E18" = 100
    "F18" = 18
    "G18" = 9
    "H18" = 100
    "I18" = 112
    "J18" = 78
    "K18" = 43.589743589743
    "L18" = -15.151515151515
    "M18" = 124
    "N18" = -73.396674584323
    "C19" = 12
    "D19" = 12
    "E19" = 0
    "F19" = 52
    "G19" = 34
    "H19" = 52.941176470588
    "I19" = 233
    "J19" = 179
    "K19" = 30.167597765363
    "L19" = 44.720496894409
    "M19" = 167.816091954023
    "N19" = 75.187969924812
    "C20" = 4
    "E20" = -20
    "G20" = 30
    "H20" = -26.666666666666
    "I20" = 99
    "J20" = 183
    "K20" = -45.901639344262
    "L20" = -12.389380530973
    "M20" = 110.63829787234
    "N20" = -55
    "C21" = 42
    "D21" = 41
    "E21" = 2.439024390243
    "F21" = 199
    "G21" = 156
    "H21" = 27.564102564102
    "I21" = 926
    "J21" = 854
    "K21" = 8.430913348946
    "L21" = 13.899138991389
    "M21" = 129.207920792079
    "N21" = -43.708206686930
    "C23" = 8
    "D23" = 6
    "E23" = 33.333333333333
    "G23" = 18
    "H23" = 50
    "I23" = 154
    "J23" = 160
    "K23" = -3.75
    "L23" = 7.692307692307
    "M23" = 105.333333333333
    "C24" = 21
    "D24" = 26
    "E24" = -19.230769230769
    "F24" = 75
    "G24" = 102
    "H24" = -26.470588235294
    "I24" = 445
    "J24" = 465
    "K24" = -4.301075268817
    "L24" = 0.678733031674
    "M24" = 51.360544217687
    "C25" = 5
    "D25" = 7
    "E25" = -28.571428571428
    "F25" = 10
    "G25" = 24
    "H25" = -58.333333333333
    "I25" = 65
    "J25" = 100
    "K25" = -35
    "L25" = -50.757575757575
    "C26" = 13
    "D26" = 26
    "E26" = -50
    "F26" = 59
    "G26" = 92
    "H26" = -35.869565217391
    "I26" = 337
    "J26" = 463
    "K26" = -27.213822894168
    "L26" = -16.169154228855
    "M26" = 0.898203592814
    "C27" = 2
    "D27" = 1
    "E27" = 100
    "F27" = 5
    "G27" = 4
    "H27" = 25
    "I27" = 26
    "J27" = 26
    "K27" = 0
    "L27" = -3.703703703703
    "C28" = 1
    "D28" = 1
    "E28" = 0
    "F28" = 11
    "G28" = 8
    "H28" = 37.5
    "I28" = 35
    "J28" = 42
    "K28" = -16.666666666666
    "L28" = 29.629629629629
    "C29" = 1
    "I29" = 10
    "K29" = 0
    "L29" = -28.571428571428
    "M29" = -52.380952380952
    "N29" = -67.741935483871
    "C30" = 1
    "I30" = 10
    "K30" = 0
    "L30" = -23.076923076923
    "M30" = -41.176470588235
    "N30" = -67.741935483871
}
foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
